# Word edit: the ".klb" directory sentence now also mentions ".tif" as an
# alternative directory extension, and the "_GoBack" bookmark (which used to
# sit at the very end of the document, after "before Export. ") moves to sit
# right before the word "directory." in that sentence.

$d = $word.ActiveDocument

$rsq = [char]0x2019   # ' (right single quotation mark)
$lsq = [char]0x2018   # ' (left single quotation mark)

# --- 1. Remove the old "_GoBack" bookmark (currently near "before Export. ") ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Find the (unique) "directory." that follows the '.klb' mention ---
$rng = $d.Content
$found = $rng.Find.Execute("directory.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insPoint = $rng.Start

    # Insert "or '.tif' " right before "directory."
    $p1 = $d.Range($insPoint, $insPoint)
    $p1.InsertBefore("or " + $lsq + ".")

    $p2 = $d.Range($p1.End, $p1.End)
    $p2.InsertBefore("tif")

    $p3 = $d.Range($p2.End, $p2.End)
    $p3.InsertBefore($rsq + " ")

    # --- 3. Re-create "_GoBack" right before "directory." ---
    $goBackRange = $d.Range($p3.End, $p3.End)
    $d.Bookmarks.Add("_GoBack", $goBackRange)
}
